$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(90553, "Ravi Lucca Camargo", "Vendas", "Consulta medica", 6, 45078, 5439.93),
    @(70346, "Ana Almeida", "Atendimento ao Cliente", "Doenca", 8, 45102, 6776.54),
    @(25863, "Cecília Teixeira", "Operacoes", "Problemas pessoais", 7, 45103, 4701.77),
    @(90534, "Mariane Campos", "Financeiro", "Problemas pessoais", 4, 45100, 6991.36),
    @(48773, "Ana Vitória da Paz", "Juridico", "Outros", 3, 45085, 9155.66),
    @(89336, "Stephany Mendonça", "Vendas", "Outros", 7, 45103, 7116.65),
    @(84183, "Maysa Moreira", "Operacoes", "Doenca", 6, 45091, 8482.360000000001),
    @(14754, "Thiago Cunha", "Financeiro", "Problemas pessoais", 7, 45087, 4143.18),
    @(97835, "Igor Nunes", "Marketing", "Outros", 4, 45101, 6133.83),
    @(49800, "Emanuel Santos", "Vendas", "Problemas pessoais", 7, 45096, 4338.84)
)

$row = 2
foreach ($rec in $data) {
    $ws.Cells.Item($row, 1).Value = $rec[0]
    $ws.Cells.Item($row, 2).Value = $rec[1]
    $ws.Cells.Item($row, 3).Value = $rec[2]
    $ws.Cells.Item($row, 4).Value = $rec[3]
    $ws.Cells.Item($row, 5).Value = $rec[4]
    $ws.Cells.Item($row, 6).Value = $rec[5]
    $ws.Cells.Item($row, 7).Value = $rec[6]
    $row++
}
